$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 236
$ws.Range("I33").Value = 134.2
$ws.Range("K33").Value = 134.2
$ws.Range("M33").Value = 94.80000000000001

$ws.Range("H64").Value = 3394.4614
$ws.Range("I64").Value = 3343.75
$ws.Range("K64").Value = 3343.75
$ws.Range("M64").Value = -3095.75

$ws.Range("H67").Value = 3394.4614
$ws.Range("I67").Value = 3343.75
$ws.Range("K67").Value = 3343.75
$ws.Range("M67").Value = -2485.75

$ws.Range("H74").Value = 4908.7144
$ws.Range("I74").Value = 5124.8823
$ws.Range("K74").Value = 5124.8823
$ws.Range("M74").Value = -4188.8823

$ws.Range("H77").Value = 4908.7144
$ws.Range("I77").Value = 5124.8823
$ws.Range("K77").Value = 25624.4115
$ws.Range("M77").Value = -20944.4115

$ws.Range("H100").Value = 3824.6667
$ws.Range("I100").Value = 2020.8462
$ws.Range("K100").Value = 2020.8462
$ws.Range("M100").Value = -1479.8462

$ws.Range("H116").Value = 4131.9
$ws.Range("I116").Value = 3825.7222
$ws.Range("K116").Value = 3825.7222
$ws.Range("M116").Value = -383.7222000000002

$ws.Range("H132").Value = 4753.5186
$ws.Range("I132").Value = 3834.2273
$ws.Range("K132").Value = 11502.6819
$ws.Range("M132").Value = -8972.6819

$ws.Range("H137").Value = 60970.844
$ws.Range("I137").Value = 112099.6
$ws.Range("J137").Value = 4161.1113
$ws.Range("K137").Value = 336298.8
$ws.Range("L137").Value = 12483.3339
$ws.Range("M137").Value = -333748.8
$ws.Range("N137").Value = -17583.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4508.5864
$ws.Range("I32").Value = 2236.4902
$ws.Range("K32").Value = 2236.4902
$ws.Range("M32").Value = -1949.4902

$ws.Range("H132").Value = 3434.8215
$ws.Range("I132").Value = 3559.6086
$ws.Range("K132").Value = 10678.8258
$ws.Range("M132").Value = -8148.825800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 52829.668
$ws.Range("J74").Value = 52829.668
$ws.Range("L74").Value = 52829.668
$ws.Range("N74").Value = -54701.668

$ws.Range("H77").Value = 52829.668
$ws.Range("J77").Value = 52829.668
$ws.Range("L77").Value = 158489.004
$ws.Range("N77").Value = -167849.004

$ws.Range("H86").Value = 1968.4375
$ws.Range("I86").Value = 1711
$ws.Range("J86").Value = 2397.5
$ws.Range("K86").Value = 1711
$ws.Range("L86").Value = 2397.5
$ws.Range("M86").Value = -588
$ws.Range("N86").Value = -4643.5

$ws.Range("H89").Value = 1968.4375
$ws.Range("I89").Value = 1711
$ws.Range("J89").Value = 2397.5
$ws.Range("K89").Value = 8555
$ws.Range("L89").Value = 11987.5
$ws.Range("M89").Value = -2939
$ws.Range("N89").Value = -23219.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 273439.6
$ws.Range("J31").Value = 4020.2778
$ws.Range("L31").Value = 4020.2778
$ws.Range("N31").Value = -4610.2778

$ws.Range("H34").Value = 273439.6
$ws.Range("J34").Value = 4020.2778
$ws.Range("L34").Value = 4020.2778
$ws.Range("N34").Value = -4424.2778

$ws.Range("H62").Value = 2415.7856
$ws.Range("J62").Value = 2197
$ws.Range("L62").Value = 2197
$ws.Range("N62").Value = -3445

$ws.Range("H65").Value = 2415.7856
$ws.Range("J65").Value = 2197
$ws.Range("L65").Value = 10985
$ws.Range("N65").Value = -17225

$ws.Range("H86").Value = 774022.3
$ws.Range("I86").Value = 1255074.5
$ws.Range("K86").Value = 1255074.5
$ws.Range("M86").Value = -1253951.5

$ws.Range("H89").Value = 774022.3
$ws.Range("I89").Value = 1255074.5
$ws.Range("K89").Value = 6275372.5
$ws.Range("M89").Value = -6269756.5

$ws.Range("H92").Value = 39750
$ws.Range("J92").Value = 39750
$ws.Range("L92").Value = 39750
$ws.Range("N92").Value = -44742

$ws.Range("H132").Value = 4123.619
$ws.Range("I132").Value = 2511.7693
$ws.Range("K132").Value = 7535.3079
$ws.Range("M132").Value = -5005.3079

$ws.Range("H134").Value = 3159.0417
$ws.Range("I134").Value = 2977.0476
$ws.Range("J134").Value = 4433
$ws.Range("K134").Value = 8931.1428
$ws.Range("L134").Value = 13299
$ws.Range("M134").Value = -6396.1428
$ws.Range("N134").Value = -18369

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4360
$ws.Range("I11").Value = 2540
$ws.Range("K11").Value = 7620
$ws.Range("M11").Value = -7480

$ws.Range("H39").Value = 951.0769
$ws.Range("J39").Value = 1653.6666
$ws.Range("L39").Value = 4960.9998
$ws.Range("N39").Value = -5548.9998

$ws.Range("H55").Value = 7088.2144
$ws.Range("J55").Value = 7629.615
$ws.Range("L55").Value = 22888.845
$ws.Range("N55").Value = -23242.845

$ws.Range("H69").Value = 5700
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622

$ws.Range("H72").Value = 5700
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3033.875
$ws.Range("I70").Value = 3031.4546
$ws.Range("K70").Value = 3031.4546
$ws.Range("M70").Value = -2761.4546

$ws.Range("H73").Value = 3033.875
$ws.Range("I73").Value = 3031.4546
$ws.Range("K73").Value = 3031.4546
$ws.Range("M73").Value = -2095.4546

$ws.Range("H123").Value = 39595.6
$ws.Range("J123").Value = 39595.6
$ws.Range("L123").Value = 39595.6
$ws.Range("N123").Value = -44495.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4411.154
$ws.Range("I46").Value = 3375.111
$ws.Range("K46").Value = 3375.111
$ws.Range("M46").Value = -3187.111

$ws.Range("H68").Value = 9365.179
$ws.Range("I68").Value = 8047.115
$ws.Range("J68").Value = 26500
$ws.Range("K68").Value = 8047.115
$ws.Range("L68").Value = 26500
$ws.Range("M68").Value = -7298.115
$ws.Range("N68").Value = -27998

$ws.Range("H71").Value = 9365.179
$ws.Range("I71").Value = 8047.115
$ws.Range("J71").Value = 26500
$ws.Range("K71").Value = 40235.575
$ws.Range("L71").Value = 132500
$ws.Range("M71").Value = -36491.575
$ws.Range("N71").Value = -139988

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
